$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New 4th column header (D1 = "C/A"), copying C1's header formatting
# (bold font + border + centered alignment) onto the new cell.
$ws.Range("D1").Value = "C/A"
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)  # xlPasteFormats

# Row labels (column A)
$ws.Range("A2").Value = "LF Lag"
$ws.Range("A3").Value = "FFR Lag"
$ws.Range("A4").Value = "C/A Lag"
$ws.Range("A5").Value = "r2"

# Row 2 (LF Lag) coefficients
$ws.Range("B2").Value = "0.818***"
$ws.Range("C2").Value = "1.303***"
$ws.Range("D2").Value = "0.875*"

# Row 3 (FFR Lag) coefficients
$ws.Range("B3").Value = "0.384**"
$ws.Range("C3").Value = "0.673***"
$ws.Range("D3").Value = "0.664**"

# Row 4 (C/A Lag) coefficients -- these look like plain numbers but must stay
# text, so force a Text number format before assigning, then restore the
# default (unstyled) formatting of a plain data cell so it keeps matching
# its siblings.
$ws.Range("B4:C4").NumberFormat = "@"
$ws.Range("B4").Value = "0.024"
$ws.Range("C4").Value = "-0.074"
$ws.Range("B5").Copy()
$ws.Range("B4:C4").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("D4").Value = "-0.48*"

# Row 5 (r2) numeric values
$ws.Range("B5").Value = 0.5558229018339893
$ws.Range("C5").Value = 0.6596486352855867
$ws.Range("D5").Value = 0.4513326749396927
